# Update the cryptocurrency price/volume snapshot (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain text (e.g. "57.456.21" with thousands-dots), not
# numbers. Some new values look like a single plain decimal (e.g. "523.92") which
# Excel would otherwise auto-convert to a Number on assignment - force those
# particular cells to keep a Text format so they round-trip as strings, matching
# the source data's inline-string cell type.
$textPriceCells = "D5","D6","D9","D12","D15","D20","D21","D22","D24","D25","D31","D32","D34","D37","D40","D42","D44","D45","D48","D49","D51"
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormatLocal = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "57.456.21"
$ws.Range("E2").Value = "  -0.72%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.107.74"
$ws.Range("E3").Value = "  +1.27%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "523.92"
$ws.Range("E5").Value = "  +1.30%  "

# Row 6 - Solana
$ws.Range("D6").Value = "141.38"
$ws.Range("E6").Value = "  -0.63%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.107.20"
$ws.Range("E8").Value = "  +1.33%  "

# Row 9 - XRP
$ws.Range("D9").Value = "0.437"
$ws.Range("E9").Value = "  +0.10%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -1.02%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +0.67%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.386"
$ws.Range("E12").Value = "  +1.67%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.642.15"
$ws.Range("E13").Value = "  +1.37%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +1.12%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "26.15"
$ws.Range("E15").Value = "  -0.17%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +0.08%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "57.533.68"
$ws.Range("E17").Value = "  -0.58%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.107.61"
$ws.Range("E18").Value = "  +1.53%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +0.68%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "12.80"
$ws.Range("E20").Value = "  -0.56%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "8.07"
$ws.Range("E21").Value = "  -0.66%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "336.96"
$ws.Range("E22").Value = "  +1.61%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.08%  "

# Row 24 - Polygon
$ws.Range("D24").Value = "0.513"
$ws.Range("E24").Value = "  +2.74%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "66.59"
$ws.Range("E25").Value = "  +1.29%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  -0.58%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("E27").Value = "  +0.28%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "0.0₃0918"
$ws.Range("E28").Value = "  +1.45%  "

# Row 29 - RenderToken
$ws.Range("E29").Value = "  +1.68%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("E30").Value = "  -0.14%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "1.86"
$ws.Range("E31").Value = "  +2.19%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "21.01"
$ws.Range("E32").Value = "  +1.42%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  +0.42%  "

# Row 34 - Monero
$ws.Range("D34").Value = "157.73"
$ws.Range("E34").Value = "  +2.05%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  +2.48%  "

# Row 36 - Aptos
$ws.Range("E36").Value = "  +2.48%  "

# Row 37 - EnergySwap
$ws.Range("D37").Value = "27.02"
$ws.Range("E37").Value = "  -0.67%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  +1.01%  "

# Row 39 - Hedera
$ws.Range("E39").Value = "  -1.78%  "

# Rows 40-42 - coins rotated (Filecoin, RenzoRestakedETH, Mantle order) with refreshed figures
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "3.96"
$ws.Range("E40").Value = "  +1.05%  "

$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "3.148.29"
$ws.Range("E41").Value = "  +1.28%  "

$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.687"
$ws.Range("E42").Value = "  +4.63%  "

# Row 43 - Stacks
$ws.Range("E43").Value = "  +10.48%  "

# Row 44 - OKB
$ws.Range("D44").Value = "36.82"
$ws.Range("E44").Value = "  +0.57%  "

# Row 45 - FirstDigitalUSD
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.01%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.307.34"
$ws.Range("E46").Value = "  +2.27%  "

# Row 47 - VeChain
$ws.Range("E47").Value = "  +0.49%  "

# Row 48 - ONDO
$ws.Range("D48").Value = "0.979"
$ws.Range("E48").Value = "  +4.00%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").Value = "20.73"
$ws.Range("E49").Value = "  +0.05%  "

# Row 50 - Cosmos
$ws.Range("E50").Value = "  +2.11%  "

# Row 51 - SuiNetwork
$ws.Range("D51").Value = "0.731"
$ws.Range("E51").Value = "  -0.22%  "
